$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 20 ---
$ws.Range("A20").Value = "ENW000010"
$ws.Range("B7").Copy() | Out-Null
$ws.Range("B20").PasteSpecial(-4122) | Out-Null
$ws.Range("B20").Value = "OPQA-1701"

$ws.Range("C7").Copy() | Out-Null
$ws.Range("C20").PasteSpecial(-4122) | Out-Null
$ws.Range("C20").Value = "Verify that the ""Thanks for your interest in EndNote......"" modal displayed when user clicks on the export button when user is signed to facebook account and not having existing steam account"

$ws.Range("D20").Value = "Y"
$ws.Range("D20").Borders.Item(7).LineStyle = 1
$ws.Range("D20").Borders.Item(10).LineStyle = 1

$ws.Rows("20").RowHeight = 45

# --- Row 21 ---
$ws.Range("A21").Value = "ENW000012"
$ws.Range("B7").Copy() | Out-Null
$ws.Range("B21").PasteSpecial(-4122) | Out-Null
$ws.Range("B21").Value = "OPQA-1701"

$ws.Range("C7").Copy() | Out-Null
$ws.Range("C21").PasteSpecial(-4122) | Out-Null
$ws.Range("C21").Value = "Verify that the ""Thanks for your interest in EndNote......"" modal displayed when user clicks on the export button when user is signed to facebook account and not having existing steam account"

$ws.Range("D21").Value = "Y"
$ws.Range("D21").Borders.Item(7).LineStyle = 1
$ws.Range("D21").Borders.Item(10).LineStyle = 1

$ws.Rows("21").RowHeight = 45

# --- Row 22 ---
$ws.Range("A22").Value = "ENW000011"
$ws.Range("B7").Copy() | Out-Null
$ws.Range("B22").PasteSpecial(-4122) | Out-Null
$ws.Range("B22").Value = "OPQA-1701"

$ws.Range("C22").Value = "steam Login  account page"
$ws.Range("C22").Borders.Item(7).LineStyle = 1
$ws.Range("C22").Borders.Item(10).LineStyle = 1
$ws.Range("C22").WrapText = $true

$ws.Range("D22").Value = "Y"
$ws.Range("D22").Borders.Item(7).LineStyle = 1
$ws.Range("D22").Borders.Item(10).LineStyle = 1

# --- Row 23 ---
$ws.Range("A23").Value = "ENW00029"
$ws.Range("B7").Copy() | Out-Null
$ws.Range("B23").PasteSpecial(-4122) | Out-Null
$ws.Range("B23").Value = "OPQA-1701"

$ws.Range("C23").Value = "Social Login account page"
$ws.Range("C23").Borders.Item(7).LineStyle = 1
$ws.Range("C23").Borders.Item(10).LineStyle = 1
$ws.Range("C23").WrapText = $true

$ws.Range("D23").Value = "Y"
$ws.Range("D23").Borders.Item(7).LineStyle = 1
$ws.Range("D23").Borders.Item(10).LineStyle = 1

# --- selection / view state ---
$ws.Range("A20:D23").Select()
